$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 19:52"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 726856
$ws.Range("C4").Value = 17121
$ws.Range("E4").Value = 624720

# Alemania (row 8)
$ws.Range("B8").Value = 142751
$ws.Range("C8").Value = 1354
$ws.Range("E8").Value = 52939
$ws.Range("G8").Value = 60
$ws.Range("H8").Value = 4412

# Brasil (row 15)
$ws.Range("B15").Value = 35025
$ws.Range("C15").Value = 1343
$ws.Range("E15").Value = 18796
$ws.Range("G15").Value = 62
$ws.Range("H15").Value = 2203

# Canada (row 16)
$ws.Range("B16").Value = 33137
$ws.Range("C16").Value = 1210
$ws.Range("D16").Value = 11126
$ws.Range("E16").Value = 20665

# Peru overtakes Suecia in total cases -> rows 23/24 swap countries
$ws.Range("A23").Value = "Peru"
$ws.Range("B23").Value = 14420
$ws.Range("C23").Value = 931
$ws.Range("D23").Value = 6541
$ws.Range("E23").Value = 7531
$ws.Range("F23").Value = 137
$ws.Range("G23").Value = 48
$ws.Range("H23").Value = 348

$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 13822
$ws.Range("C24").Value = 606
$ws.Range("D24").Value = 550
$ws.Range("E24").Value = 11761
$ws.Range("F24").Value = 1054
$ws.Range("G24").Value = 111
$ws.Range("H24").Value = 1511

# Israel (row 25)
$ws.Range("B25").Value = 13265
$ws.Range("C25").Value = 283
$ws.Range("D25").Value = 3456
$ws.Range("E25").Value = 9645
$ws.Range("F25").Value = 164
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 164

# Ecuador overtakes Polonia in total cases -> rows 29/30 swap countries
$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 9022
$ws.Range("C29").Value = 572
$ws.Range("D29").Value = 1008
$ws.Range("E29").Value = 7558
$ws.Range("F29").Value = 168
$ws.Range("G29").Value = 35
$ws.Range("H29").Value = 456

$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 8742
$ws.Range("C30").Value = 363
$ws.Range("D30").Value = 981
$ws.Range("E30").Value = 7414
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = 347

# Bielorrusia (row 47)
$ws.Range("E47").Value = 4392
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 45

# Egipto (row 53)
$ws.Range("B53").Value = 3032
$ws.Range("C53").Value = 188
$ws.Range("D53").Value = 701
$ws.Range("E53").Value = 2107
$ws.Range("G53").Value = 19
$ws.Range("H53").Value = 224

# Republica de Chipre (row 87)
$ws.Range("D87").Value = 79
$ws.Range("E87").Value = 670

# Jordania overtakes Reunion in total cases -> rows 107/108 swap countries
$ws.Range("A107").Value = "Jordania"
$ws.Range("B107").Value = 413
$ws.Range("C107").Value = 6
$ws.Range("D107").Value = 269
$ws.Range("E107").Value = 137
$ws.Range("F107").Value = 5
$ws.Range("H107").Value = 7

$ws.Range("A108").Value = "Reunion"
$ws.Range("B108").Value = 407
$ws.Range("C108").Value = 5
$ws.Range("D108").Value = 237
$ws.Range("E108").Value = 170
$ws.Range("F108").Value = 4
$ws.Range("H108").Value = 0

# Madagascar (row 136)
$ws.Range("B136").Value = 120
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 35
$ws.Range("E136").Value = 85

# Aruba (row 141)
$ws.Range("D141").Value = 44
$ws.Range("E141").Value = 50

# Curazao (row 192)
$ws.Range("D192").Value = 11
$ws.Range("E192").Value = 2

# Nicaragua (row 202)
$ws.Range("E202").Value = 1
$ws.Range("G202").Value = 1
$ws.Range("H202").Value = 2
